$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 42/43 (entered first, fixes shared-string ordering) ---
$ws.Range("C43").Value = "6 to 11_02"
$ws.Range("C42").Value = "1 to 11_01"

# --- Header row (row 2) ---
# D2, F2, G2 stay the same (SimilarText / Wrong / CorrectText)
$ws.Range("I2").Value = "SumWrong"
$ws.Range("J2").Value = "SumCorrectText"
$ws.Range("K2").Value = "Accuracy"
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()

# --- Row 3 (data row 1) ---
$ws.Range("C3").Value = "1_1 to 11_01"
$ws.Range("D3").Value = 2297
$ws.Range("F3").Value = 101
$ws.Range("G3").Value = 2347
$ws.Range("I3").Formula = "=SUM(F3:F100)"
$ws.Range("J3").Formula = "=SUM(G3:G100)"
$ws.Range("K3").Formula = "=1-(I3/J3)"
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()

# --- Row 8 (data row 6) ---
$ws.Range("C8").Value = "6_1 to 11_02"
$ws.Range("D8").Value = 2576
$ws.Range("F8").Value = 50
$ws.Range("G8").Value = 2620

# --- Column widths ---
# (target stored widths are 10.3984375 / 14 / 14.296875; the host rounds
#  ColumnWidth to the nearest 1/7 "pixel" unit on write, so these inputs are
#  chosen to land on the closest representable stored width)
$ws.Columns.Item(9).ColumnWidth = 9.714285714285714
$ws.Columns.Item(10).ColumnWidth = 13.285714285714286
$ws.Columns.Item(11).ColumnWidth = 13.571428571428571

# --- Sheet view: remove frozen/topLeftCell override, update selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("L29").Select()
